$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 3 (Sierra Nevada entry) entirely - the sheet now has
# only a header row plus a single data row.
$ws.Rows(3).Delete()

# Update row 2: column A now holds the (hyperlink-styled) source URL,
# column B the contract date, column C the contract description.
$ws.Range("A2").Value = "https://www.defense.gov/News/Contracts/Contract/Article/2644648/"
$ws.Range("B2").Value = 44350
$ws.Range("C2").Value = "Northrop Grumman Systems Corp., Aerospace Systems, Melbourne, Florida, is awarded a `$12,015,026 modification (P00036) to a previously awarded cost-plus-fixed-fee contract (N0001914C0036). This modification increases the ceiling to extend services and adds hours increasing the full-scale fatigue repair time to achieve the required simulated flight hours in support of E-2D Advanced Hawkeye aircraft development. Work will be performed in El Segundo, California (59%); Melbourne, Florida (35%); and Bethpage, New York (6%), and is expected to be completed in June 2023. No funds will be obligated at time of award. The Naval Air Systems Command, Patuxent River, Maryland, is the contracting activity."

# Give the URL cell the built-in "Hyperlink" look (underline + theme color)
# while keeping its own font name/size.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 12

# Wrap text on both the date and description cells, row grew tall to fit.
$ws.Range("A2").WrapText = $true
$ws.Range("B2").WrapText = $true
$ws.Rows(2).RowHeight = 90

$ws.Range("A2").Select()
